$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at I, shifting the old "testm5.brms" / "no conifers"
# column (old I) one slot to the right (new J). The new, empty I column is
# where the new "test5cent.rstanarm" model results go.
$ws.Columns("I").Insert()

# --- New column I: test5cent.rstanarm (no conifers) -----------------------
$ws.Range("I1").Value = "no conifers"
$ws.Range("I2").Value = "test5cent.rstanarm"
$ws.Range("I3").Value = "bb"
$ws.Range("I4").Value = 47
$ws.Range("I5").Value = 5
$ws.Range("I6").Value = 9514
$ws.Range("I7").Value = 96.169168999999997
$ws.Range("I8").Value = -9.3518709999999992
$ws.Range("I9").Value = -1.6437349999999999
$ws.Range("I10").Value = 0.72115099999999999
$ws.Range("I11").Value = "4 divergent transitions"

# --- Fill in newly-populated cells in column H -----------------------------
$ws.Range("H4").Value = 47
$ws.Range("H7").Value = 95.28
$ws.Range("H8").Value = -9.69
$ws.Range("H9").Value = -1.69
$ws.Range("H10").Value = 0.7

# --- Remove the old scratch prediction block (rows 21-23, cols F:H) -------
$ws.Range("F21:H23").ClearContents()
